$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "Pasar a siguiente ronda"
$ws.Range("B17").Value = "15/ID_partida"
$ws.Range("C17").Value = "-"
$ws.Range("D17").Value = "15`$ID_partida/indexJugador/fichas/puntos/numCartas/cartas"
$ws.Range("E17").Value = "Cuando termina una ronda, el servidor envia notificaciones a los jugadores con las cartas, puntos y fichas de cada jugador y decide quien ha ganado la ronda. El cliente observa los resultados y responde pidiendo que comienze la siguiente."

$ws.Rows.Item(17).RowHeight = 45

$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("E18").Select()
